$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently spans A1:D15 with headers:
#   id_building_efficiency_class | id_radiator | unit | value
# We need to:
#   - rename the "value" column to "space_heating"
#   - add a new "hot_water" column (E)
#   - fill in supply-temperature values for both columns, per row

$tbl = $ws.ListObjects.Item(1)

# Grow the table to include the new column E before touching headers,
# otherwise the header text we set below gets overwritten by a default name.
$tbl.Resize($ws.Range("A1:E15"))

# Update header row
$ws.Range("D1").Value = "space_heating"
$ws.Range("E1").Value = "hot_water"

# Supply temperature (degree) values:
#   id_radiator = 1 (rows 2-8)  -> space_heating = 45
#   id_radiator = 2 (rows 9-15) -> space_heating = 35
#   hot_water is always 60 regardless of radiator type
$spaceHeatingValues = @(45, 45, 45, 45, 45, 45, 45, 35, 35, 35, 35, 35, 35, 35)

for ($i = 0; $i -lt $spaceHeatingValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $spaceHeatingValues[$i]
    $ws.Cells.Item($row, 5).Value = 60
}

# Best-fit the new columns' widths like Excel would after entering data
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# Match the final selection shown in the saved workbook
$ws.Range("D12").Select() | Out-Null
